# Add the missing "EndophilinA1" amphipathic-helix sequence as a new row
# at the bottom of the data table (row 42).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Protein name is entered before the sequence so that shared-string indices
# line up the same way Excel's own "type it in, tab across" flow would.
$ws.Range("B42").Value = "EndophilinA1"
$ws.Range("A42").Value = "SVAGLKKQFHKATQKVSEKV"
$ws.Range("C42").Value = 1

# Leave the view scrolled down with the new row's first cell selected,
# mirroring where the user would naturally end up after typing the row.
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("A42").Select()
